$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1376
$ws.Range("F5").Value = 5713
$ws.Range("F6").Value = 443
$ws.Range("F7").Value = 1044
$ws.Range("F8").Value = 3330
$ws.Range("F9").Value = 6522
$ws.Range("F10").Value = 192
$ws.Range("F11").Value = 1271
$ws.Range("F12").Value = 736
$ws.Range("F13").Value = 89
$ws.Range("F14").Value = 7
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 1106
$ws.Range("F18").Value = 87
$ws.Range("F20").Value = 158
$ws.Range("F22").Value = 929
$ws.Range("G22").Value = 89
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 88
$ws.Range("F28").Value = 1137
$ws.Range("F29").Value = 10
$ws.Range("F30").Value = 15
$ws.Range("F33").Value = 252
$ws.Range("F34").Value = 45
$ws.Range("F35").Value = 252
$ws.Range("F36").Value = 1156
$ws.Range("F37").Value = 49
$ws.Range("F38").Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F24").Value = 601
$ws.Range("F27").Value = 108
$ws.Range("F28").Value = 640
$ws.Range("F29").Value = 945
$ws.Range("F32").Value = 76

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 530
$ws.Range("F7").Value = 279
$ws.Range("F8").Value = 582

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1376
$ws.Range("F8").Value = 530
$ws.Range("F9").Value = 279
$ws.Range("F10").Value = 279
$ws.Range("F13").Value = 5713
$ws.Range("F14").Value = 443
$ws.Range("F15").Value = 1044
$ws.Range("F16").Value = 3330
$ws.Range("F18").Value = 6522
$ws.Range("F19").Value = 192
$ws.Range("F20").Value = 1271
$ws.Range("F24").Value = 736
$ws.Range("F25").Value = 89
$ws.Range("F26").Value = 1106
$ws.Range("F28").Value = 87
$ws.Range("F29").Value = 158
$ws.Range("F31").Value = 929
$ws.Range("G31").Value = 89
$ws.Range("F32").Value = 601
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 88
$ws.Range("F35").Value = 1137
$ws.Range("F38").Value = 945
$ws.Range("F41").Value = 252
$ws.Range("F42").Value = 45
$ws.Range("F43").Value = 76
$ws.Range("F44").Value = 252
$ws.Range("F49").Value = 80
